# Add the new Low/Hi/Hi_Hi alarm limit columns to the DAS Facts header row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (K1:M1), continuing the existing header row.
$ws.Range("K1").Value = "Low Limit"
$ws.Range("L1").Value = "Hi Limit"
$ws.Range("M1").Value = "Hi_Hi Limit"

# Match the formatting of the existing header cells (e.g. J1 "Data Source")
# by copying its format onto the new header cells.
$ws.Range("J1").Copy() | Out-Null
$ws.Range("K1:M1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Update the active selection to reflect where the author left off editing.
$ws.Range("K7").Select() | Out-Null
